# Updates the cryptos list with refreshed price/volume data
# (commit: "Updated cryptos list on Fri Jun 21 18:41:55 UTC 2024 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.914.37'
$ws.Range("E2").Value = '  -1.74%  '
$ws.Range("D3").Value = '3.495.34'
$ws.Range("E3").Value = '  -0.94%  '
$ws.Range("D4").Value = "'" + '1.00'
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = "'" + '583.84'
$ws.Range("E5").Value = '  -1.85%  '
$ws.Range("D6").Value = "'" + '130.65'
$ws.Range("E6").Value = '  -2.95%  '
$ws.Range("D7").Value = '3.493.57'
$ws.Range("E7").Value = '  -0.88%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E9").Value = '  -2.04%  '
$ws.Range("D10").Value = "'" + '0.124'
$ws.Range("E10").Value = '  -0.83%  '
$ws.Range("D11").Value = "'" + '7.11'
$ws.Range("E11").Value = '  -0.15%  '
$ws.Range("D12").Value = "'" + '0.376'
$ws.Range("E12").Value = '  -3.08%  '
$ws.Range("D13").Value = '4.079.56'
$ws.Range("E13").Value = '  -1.27%  '
$ws.Range("D14").Value = "'" + '27.28'
$ws.Range("E14").Value = '  -1.51%  '
$ws.Range("E15").Value = '  +1.10%  '
$ws.Range("D16").Value = '3.491.75'
$ws.Range("E16").Value = '  -1.19%  '
$ws.Range("D17").Value = "'" + '0.0000177'
$ws.Range("E17").Value = '  -2.93%  '
$ws.Range("D18").Value = '63.981.90'
$ws.Range("E18").Value = '  -1.63%  '
$ws.Range("E19").Value = '  -3.43%  '
$ws.Range("D20").Value = "'" + '14.09'
$ws.Range("E20").Value = '  -2.50%  '
$ws.Range("D21").Value = "'" + '5.61'
$ws.Range("E21").Value = '  -1.81%  '
$ws.Range("D22").Value = "'" + '380.53'
$ws.Range("E22").Value = '  -2.94%  '
$ws.Range("D23").Value = "'" + '0.569'
$ws.Range("E23").Value = '  -1.95%  '
$ws.Range("D24").Value = '3.629.19'
$ws.Range("E24").Value = '  -1.25%  '
$ws.Range("D25").Value = "'" + '73.21'
$ws.Range("E25").Value = '  -1.89%  '
$ws.Range("E26").Value = '  +0.16%  '
$ws.Range("E27").Value = '  +1.61%  '
$ws.Range("D28").Value = "'" + '1.55'
$ws.Range("E28").Value = '  -2.61%  '
$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D29").Value = "'" + '7.45'
$ws.Range("E29").Value = '  -4.13%  '
$ws.Range("B30").Value = 'Binance-PegBSC-USD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D30").Value = "'" + '0.999'
$ws.Range("E30").Value = '  -0.01%  '
$ws.Range("D31").Value = "'" + '8.22'
$ws.Range("E31").Value = '  -1.94%  '
$ws.Range("D32").Value = "'" + '2.22'
$ws.Range("E32").Value = '  -2.51%  '
$ws.Range("D33").Value = '3.500.43'
$ws.Range("E33").Value = '  -0.94%  '
$ws.Range("E34").Value = '  +0.01%  '
$ws.Range("D35").Value = "'" + '23.33'
$ws.Range("E35").Value = '  -3.72%  '
$ws.Range("D36").Value = "'" + '0.144'
$ws.Range("E36").Value = '  -0.63%  '
$ws.Range("D37").Value = "'" + '5.28'
$ws.Range("E37").Value = '  -0.16%  '
$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D38").Value = "'" + '1.55'
$ws.Range("E38").Value = '  -1.79%  '
$ws.Range("B39").Value = 'Aptos'
$ws.Range("C39").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D39").Value = "'" + '6.85'
$ws.Range("E39").Value = '  -2.10%  '
$ws.Range("D40").Value = "'" + '159.24'
$ws.Range("E40").Value = '  -5.63%  '
$ws.Range("D41").Value = "'" + '0.0788'
$ws.Range("E41").Value = '  -3.45%  '
$ws.Range("D42").Value = "'" + '0.809'
$ws.Range("E42").Value = '  -2.04%  '
$ws.Range("D43").Value = "'" + '26.04'
$ws.Range("E43").Value = '  +0.23%  '
$ws.Range("D44").Value = "'" + '1.00'
$ws.Range("E44").Value = '  +0.01%  '
$ws.Range("D45").Value = "'" + '41.87'
$ws.Range("E45").Value = '  -2.50%  '
$ws.Range("B46").Value = 'Filecoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D46").Value = "'" + '4.37'
$ws.Range("E46").Value = '  -1.59%  '
$ws.Range("B47").Value = 'ONDO'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D47").Value = "'" + '1.19'
$ws.Range("E47").Value = '  -5.82%  '
$ws.Range("E48").Value = '  -3.26%  '
$ws.Range("B49").Value = 'Cosmos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D49").Value = "'" + '6.80'
$ws.Range("E49").Value = '  -1.79%  '
$ws.Range("B50").Value = 'Maker'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D50").Value = '2.413.05'
$ws.Range("E50").Value = '  -0.28%  '
$ws.Range("D51").Value = "'" + '0.895'
$ws.Range("E51").Value = '  -1.72%  '
